# Update graphs with 19 and 20 April 2020.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 40 (19 April 2020)
$ws.Range("E40").Value = 720630
$ws.Range("E40").NumberFormat = "#,##0"
$ws.Range("F40").Value = 37202
$ws.Range("F40").NumberFormat = "#,##0"
$ws.Range("G40").Formula = "=F40-F39"
$ws.Range("G40").Style = "Normal"
$ws.Range("H40").Formula = "=F40/E40"
$ws.Range("H40").NumberFormat = "0.00%"
$ws.Range("I40").Formula = "=(F40-F39)/(E40-E39)"
$ws.Range("I40").NumberFormat = "0.00%"

# Row 41 (20 April 2020)
$ws.Range("E41").Value = 746625
$ws.Range("E41").NumberFormat = "#,##0"
$ws.Range("F41").Value = 39083
$ws.Range("F41").NumberFormat = "#,##0"
$ws.Range("G41").Formula = "=F41-F40"
$ws.Range("G41").Style = "Normal"
$ws.Range("H41").Formula = "=F41/E41"
$ws.Range("H41").NumberFormat = "0.00%"
$ws.Range("I41").Formula = "=(F41-F40)/(E41-E40)"
$ws.Range("I41").NumberFormat = "0.00%"
